$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 18 ("Lesquels sortent du lot?") - content placeholder shape
# -----------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$tr18 = $s18.Shapes.Item(2).TextFrame.TextRange

# 1) "... pour la simplicité de sa syntaxe" -> "... et ses performances"
$f = $tr18.Find(" pour la simplicité de sa syntaxe")
$f.Text = " pour la simplicité de sa syntaxe et ses performances"

# 2) "... et SQL, et pour ses fonctionnalités riches" -> "..., et ses performances"
$f = $tr18.Find(" et SQL, et pour ses fonctionnalités riches")
$f.Text = " et SQL, et pour ses fonctionnalités riches, et ses performances"

# -----------------------------------------------------------------
# Slide 19 ("Qu'en conclure ?") - content placeholder shape
# Add a new sub-bullet paragraph after "Seul le futur le dira..."
# -----------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$tr19 = $s19.Shapes.Item(2).TextFrame.TextRange

$lastPara = $tr19.Paragraphs(5)
[void]$lastPara.InsertAfter([char]13 + "EntityFramework")

$newPara = $tr19.Paragraphs(6)
[void]$newPara.InsertAfter(" est plus facile à tester unitairement en ")

$newPara = $tr19.Paragraphs(6)
[void]$newPara.InsertAfter("mockant")

$newPara = $tr19.Paragraphs(6)
[void]$newPara.InsertAfter(" son ")

$newPara = $tr19.Paragraphs(6)
[void]$newPara.InsertAfter("DbContext")

# -----------------------------------------------------------------
# Slide 20 ("À surveiller….") - content placeholder shape
# -----------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$tr20 = $s20.Shapes.Item(2).TextFrame.TextRange

$f = $tr20.Find("Xpoco")
$f.Text = "XPoco"

$f = $tr20.Find("Nreco")
$f.Text = "NReco"
